$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F6").Value = 0
$ws.Range("F8").Value = 20
$ws.Range("F9").Value = 1003
$ws.Range("F10").Value = 794
$ws.Range("F11").Value = 229
$ws.Range("F14").Value = 806
$ws.Range("F15").Value = 270
$ws.Range("F16").Value = 575
$ws.Range("F20").Value = 644
$ws.Range("F21").Value = 1153
$ws.Range("F22").Value = 2842
$ws.Range("F23").Value = 1374
$ws.Range("F24").Value = 684
$ws.Range("F28").Value = 996
$ws.Range("F29").Value = 344
$ws.Range("F30").Value = 2618
$ws.Range("F31").Value = 461
$ws.Range("F33").Value = 1378

$ws = $wb.Worksheets.Item("演出")
$ws.Range("G3").Value = "不可售"
$ws.Range("F9").Value = 40

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("G8").Value = "不可售"
$ws.Range("F12").Value = 947
$ws.Range("F15").Value = 20
$ws.Range("F16").Value = 1003
$ws.Range("F17").Value = 794
$ws.Range("F18").Value = 229
$ws.Range("F21").Value = 40
$ws.Range("F26").Value = 806
$ws.Range("F27").Value = 270
$ws.Range("F28").Value = 575
$ws.Range("F32").Value = 644
$ws.Range("F33").Value = 1153
$ws.Range("F34").Value = 2842
$ws.Range("F35").Value = 1374
$ws.Range("F36").Value = 684
$ws.Range("F42").Value = 996
$ws.Range("F43").Value = 344
$ws.Range("F44").Value = 2618
$ws.Range("F45").Value = 461
$ws.Range("F47").Value = 1378
